$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 501-502, shifting the existing data
# (rows 501:604) down to (rows 503:606). This mirrors the weekly price
# update: a new reporting date's two records (Primera / Segunda) are
# prepended to the historical series kept in this sheet.
$ws.Rows("501:502").Insert()

# New row 501 - "Primera" grade record for the new date
$ws.Range("A501").Value = 8
$ws.Range("B501").Value = "Terminal La Palmera de La Serena"
$ws.Range("C501").Value = "Coquimbo"
$ws.Range("D501").Value = 45005
$ws.Range("E501").Value = 4
$ws.Range("F501").Value = 100112017
$ws.Range("G501").Value = "Apio"
$ws.Range("H501").Value = "Americana (o)"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 1460
$ws.Range("K501").Value = 7000
$ws.Range("L501").Value = 8000
$ws.Range("M501").Value = 7500
$ws.Range("N501").Value = "`$/docena de matas"
$ws.Range("O501").Value = "Provincia del Elquí"
$ws.Range("P501").Value = 1250
$ws.Range("Q501").Value = 6
$ws.Range("R501").Value = "Hortaliza"

# New row 502 - "Segunda" grade record for the new date
$ws.Range("A502").Value = 8
$ws.Range("B502").Value = "Terminal La Palmera de La Serena"
$ws.Range("C502").Value = "Coquimbo"
$ws.Range("D502").Value = 45005
$ws.Range("E502").Value = 4
$ws.Range("F502").Value = 100112017
$ws.Range("G502").Value = "Apio"
$ws.Range("H502").Value = "Americana (o)"
$ws.Range("I502").Value = "Segunda"
$ws.Range("J502").Value = 880
$ws.Range("K502").Value = 5000
$ws.Range("L502").Value = 6000
$ws.Range("M502").Value = 5500
$ws.Range("N502").Value = "`$/docena de matas"
$ws.Range("O502").Value = "Provincia del Elquí"
$ws.Range("P502").Value = 917
$ws.Range("Q502").Value = 6
$ws.Range("R502").Value = "Hortaliza"
